$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.9039155943081556
$ws.Range("C2").Value = 0.9987679826329908
$ws.Range("D2").Value = 0.9988129942789321
$ws.Range("E2").Value = 0.9917829357034532
$ws.Range("F2").Value = 0.9964546376082618
$ws.Range("G2").Value = 0.9935453628626348
$ws.Range("H2").Value = 0.9929116505919509
$ws.Range("I2").Value = 0.8785558378176582
$ws.Range("B3").Value = 0.007054299581795931
$ws.Range("C3").Value = 0.0000772136845625937
$ws.Range("D3").Value = 0.00007579336670460179
$ws.Range("E3").Value = 0.00004824167990591377
$ws.Range("F3").Value = 0.00006705630948999897
$ws.Range("G3").Value = 0.00009689725993666798
$ws.Range("H3").Value = 0.0000008745918762542715
$ws.Range("I3").Value = 134.5385284423828
$ws.Range("B4").Value = 0.08398988097906113
$ws.Range("C4").Value = 0.008787131868302822
$ws.Range("D4").Value = 0.00870593823492527
$ws.Range("E4").Value = 0.006945623084902763
$ws.Range("F4").Value = 0.008144588209688663
$ws.Range("G4").Value = 0.009843640960752964
$ws.Range("H4").Value = 0.0009351961780339479
$ws.Range("I4").Value = 11.5990743637085
$ws.Range("B5").Value = 0.01926547475159168
$ws.Range("C5").Value = 0.002964772982522845
$ws.Range("D5").Value = 0.003025685669854283
$ws.Range("E5").Value = 0.002364977030083537
$ws.Range("F5").Value = 0.002785254968330264
$ws.Range("G5").Value = 0.003383800620213151
$ws.Range("H5").Value = 0.0003872609813697636
$ws.Range("I5").Value = 2.974501132965088
